$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.209.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.249.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.53'
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0941'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.46%  '
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.582.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.251.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.073.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0979'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("E22").Value = '  +5.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +24.30%  '
$ws.Range("E28").Value = '  -3.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("E30").Value = '  -3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0823'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.36%  '
$ws.Range("E33").Value = '  -6.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.72%  '
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("E40").Value = '  -5.21%  '
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '62.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '107.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.66%  '
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -4.02%  '
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("E51").Value = '  +15.80%  '
